$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the sheet/tab
$ws.Name = "Partner Solver Weights"

# 2. Row 25: challenge_weights/needs_weights/stage_weights correction
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 1

# 3. Row 32: weights were mistakenly stored as plain numbers; they should hold
#    the (text) weight codes "2","1","1","2" that belong to this partner.
$ws.Range("C32").Value = "'2"
$ws.Range("D32").Value = "'1"
$ws.Range("E32").Value = "'1"
$ws.Range("F32").Value = "'2"
$ws.Range("C32:F32").Style = "Normal"

# 4. Row 72: challenge_weights/needs_weights/stage_weights correction
$ws.Range("D72").Value = 3
$ws.Range("E72").Value = 2
$ws.Range("F72").Value = 4

# 5. Row 585: challenge_weights/needs_weights correction
$ws.Range("D585").Value = 1
$ws.Range("E585").Value = 1

# 6. Row 745: these cells had incorrectly been stored as text weight codes;
#    they should simply be the numeric weight 1.
$ws.Range("C745").Value = 1
$ws.Range("D745").Value = 1
$ws.Range("E745").Value = 1
$ws.Range("F745").Value = 1
